# Auto-generated edit script: updates cached market-data columns (H:N)
# on the Ixion_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to match a refreshed pull from the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# ALC!row74
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 3306
$ws.Cells.Item(74, 9).Value = 3475.375
$ws.Cells.Item(74, 11).Value = 3475.375
$ws.Cells.Item(74, 13).Value = -2539.375

# ALC!row77
$ws.Cells.Item(77, 8).Value = 3306
$ws.Cells.Item(77, 9).Value = 3475.375
$ws.Cells.Item(77, 11).Value = 17376.875
$ws.Cells.Item(77, 13).Value = -12696.875

# ALC!row92
$ws.Cells.Item(92, 8).Value = 63132092
$ws.Cells.Item(92, 9).Value = 5556357.5
$ws.Cells.Item(92, 10).Value = 111111864
$ws.Cells.Item(92, 11).Value = 5556357.5
$ws.Cells.Item(92, 12).Value = 111111864
$ws.Cells.Item(92, 13).Value = -5555109.5
$ws.Cells.Item(92, 14).Value = -111114360

# ALC!row112
$ws.Cells.Item(112, 8).Value = 30076210
$ws.Cells.Item(112, 10).Value = 35715424
$ws.Cells.Item(112, 12).Value = 107146272
$ws.Cells.Item(112, 14).Value = -107148488

# ALC!row116
$ws.Cells.Item(116, 8).Value = 21541
$ws.Cells.Item(116, 9).Value = 21541
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 21541
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = -18099
$ws.Cells.Item(116, 14).ClearContents()

# ALC!row129
$ws.Cells.Item(129, 8).Value = 987.92645
$ws.Cells.Item(129, 10).Value = 1069.9661
$ws.Cells.Item(129, 12).Value = 3209.8983
$ws.Cells.Item(129, 14).Value = -13209.8983

# ALC!row137
$ws.Cells.Item(137, 8).Value = 1684.7858
$ws.Cells.Item(137, 9).Value = 1465.625
$ws.Cells.Item(137, 10).Value = 2999.75
$ws.Cells.Item(137, 11).Value = 4396.875
$ws.Cells.Item(137, 12).Value = 8999.25
$ws.Cells.Item(137, 13).Value = -1846.875
$ws.Cells.Item(137, 14).Value = -14099.25

# ARM!row26
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 2601
$ws.Cells.Item(26, 9).Value = 1058.4286
$ws.Cells.Item(26, 10).Value = 8000
$ws.Cells.Item(26, 11).Value = 1058.4286
$ws.Cells.Item(26, 12).Value = 8000
$ws.Cells.Item(26, 13).Value = -728.4286
$ws.Cells.Item(26, 14).Value = -8660

# ARM!row32
$ws.Cells.Item(32, 8).Value = 9004.522999999999
$ws.Cells.Item(32, 9).Value = 6906.6777
$ws.Cells.Item(32, 11).Value = 6906.6777
$ws.Cells.Item(32, 13).Value = -6619.6777

# ARM!row45
$ws.Cells.Item(45, 8).Value = 10557.182
$ws.Cells.Item(45, 9).Value = 10557.182
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 10557.182
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = -10180.182
$ws.Cells.Item(45, 14).ClearContents()

# ARM!row61
$ws.Cells.Item(61, 8).Value = 4098.8545
$ws.Cells.Item(61, 9).Value = 4011.775
$ws.Cells.Item(61, 10).Value = 4331.067
$ws.Cells.Item(61, 11).Value = 4011.775
$ws.Cells.Item(61, 12).Value = 4331.067
$ws.Cells.Item(61, 13).Value = -3799.775
$ws.Cells.Item(61, 14).Value = -4755.067

# ARM!row74
$ws.Cells.Item(74, 8).Value = 1993.2963
$ws.Cells.Item(74, 9).Value = 1872.25
$ws.Cells.Item(74, 10).Value = 2090.1333
$ws.Cells.Item(74, 11).Value = 1872.25
$ws.Cells.Item(74, 12).Value = 2090.1333
$ws.Cells.Item(74, 13).Value = -998.25
$ws.Cells.Item(74, 14).Value = -3838.1333

# ARM!row77
$ws.Cells.Item(77, 8).Value = 1993.2963
$ws.Cells.Item(77, 9).Value = 1872.25
$ws.Cells.Item(77, 10).Value = 2090.1333
$ws.Cells.Item(77, 11).Value = 9361.25
$ws.Cells.Item(77, 12).Value = 10450.6665
$ws.Cells.Item(77, 13).Value = -4993.25
$ws.Cells.Item(77, 14).Value = -19186.6665

# ARM!row97
$ws.Cells.Item(97, 8).Value = 1145.4
$ws.Cells.Item(97, 9).Value = 1110.7368
$ws.Cells.Item(97, 10).Value = 1255.1666
$ws.Cells.Item(97, 11).Value = 1110.7368
$ws.Cells.Item(97, 12).Value = 1255.1666
$ws.Cells.Item(97, 13).Value = -614.7367999999999
$ws.Cells.Item(97, 14).Value = -2247.1666

# ARM!row110
$ws.Cells.Item(110, 8).Value = 1154.6364
$ws.Cells.Item(110, 9).Value = 961.95654
$ws.Cells.Item(110, 10).Value = 1597.8
$ws.Cells.Item(110, 11).Value = 961.95654
$ws.Cells.Item(110, 12).Value = 1597.8
$ws.Cells.Item(110, 13).Value = 1083.04346
$ws.Cells.Item(110, 14).Value = -5687.8

# ARM!row132
$ws.Cells.Item(132, 8).Value = 1925483.8
$ws.Cells.Item(132, 9).Value = 1494.0857
$ws.Cells.Item(132, 10).Value = 5886639
$ws.Cells.Item(132, 11).Value = 4482.257100000001
$ws.Cells.Item(132, 12).Value = 17659917
$ws.Cells.Item(132, 13).Value = -1952.257100000001
$ws.Cells.Item(132, 14).Value = -17664977

# ARM!row136
$ws.Cells.Item(136, 8).Value = 4098.8545
$ws.Cells.Item(136, 9).Value = 4011.775
$ws.Cells.Item(136, 10).Value = 4331.067
$ws.Cells.Item(136, 11).Value = 12035.325
$ws.Cells.Item(136, 12).Value = 12993.201
$ws.Cells.Item(136, 13).Value = -9485.325000000001
$ws.Cells.Item(136, 14).Value = -18093.201

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1982.9524
$ws.Cells.Item(86, 9).Value = 1777.3636
$ws.Cells.Item(86, 11).Value = 1777.3636
$ws.Cells.Item(86, 13).Value = -654.3635999999999

# BSM!row89
$ws.Cells.Item(89, 8).Value = 1982.9524
$ws.Cells.Item(89, 9).Value = 1777.3636
$ws.Cells.Item(89, 11).Value = 8886.817999999999
$ws.Cells.Item(89, 13).Value = -3270.817999999999

# BSM!row94
$ws.Cells.Item(94, 8).Value = 1533.5
$ws.Cells.Item(94, 9).Value = 822.4
$ws.Cells.Item(94, 10).Value = 2422.375
$ws.Cells.Item(94, 11).Value = 822.4
$ws.Cells.Item(94, 12).Value = 2422.375
$ws.Cells.Item(94, 13).Value = -371.4
$ws.Cells.Item(94, 14).Value = -3324.375

# BSM!row99
$ws.Cells.Item(99, 8).Value = 38462576
$ws.Cells.Item(99, 9).Value = 62500810
$ws.Cells.Item(99, 10).Value = 1398.7
$ws.Cells.Item(99, 11).Value = 62500810
$ws.Cells.Item(99, 12).Value = 1398.7
$ws.Cells.Item(99, 13).Value = -62499312
$ws.Cells.Item(99, 14).Value = -4394.7

# BSM!row105
$ws.Cells.Item(105, 8).Value = 17856.154
$ws.Cells.Item(105, 9).Value = 27753.75
$ws.Cells.Item(105, 10).Value = 2020
$ws.Cells.Item(105, 11).Value = 27753.75
$ws.Cells.Item(105, 12).Value = 2020
$ws.Cells.Item(105, 13).Value = -26006.75
$ws.Cells.Item(105, 14).Value = -5514

# BSM!row107
$ws.Cells.Item(107, 8).Value = 1115.5
$ws.Cells.Item(107, 9).Value = 1026.6111
$ws.Cells.Item(107, 11).Value = 1026.6111
$ws.Cells.Item(107, 13).Value = 893.3888999999999

# CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1925.7858
$ws.Cells.Item(16, 9).Value = 1694.4706
$ws.Cells.Item(16, 11).Value = 1694.4706
$ws.Cells.Item(16, 13).Value = -1407.4706

# CRP!row31
$ws.Cells.Item(31, 8).Value = 4842.983
$ws.Cells.Item(31, 9).Value = 1660.1666
$ws.Cells.Item(31, 10).Value = 9824.781999999999
$ws.Cells.Item(31, 11).Value = 1660.1666
$ws.Cells.Item(31, 12).Value = 9824.781999999999
$ws.Cells.Item(31, 13).Value = -1365.1666
$ws.Cells.Item(31, 14).Value = -10414.782

# CRP!row34
$ws.Cells.Item(34, 8).Value = 4842.983
$ws.Cells.Item(34, 9).Value = 1660.1666
$ws.Cells.Item(34, 10).Value = 9824.781999999999
$ws.Cells.Item(34, 11).Value = 1660.1666
$ws.Cells.Item(34, 12).Value = 9824.781999999999
$ws.Cells.Item(34, 13).Value = -1458.1666
$ws.Cells.Item(34, 14).Value = -10228.782

# CRP!row99
$ws.Cells.Item(99, 8).Value = 5501
$ws.Cells.Item(99, 9).Value = 11899.8
$ws.Cells.Item(99, 10).Value = 1501.75
$ws.Cells.Item(99, 11).Value = 11899.8
$ws.Cells.Item(99, 12).Value = 1501.75
$ws.Cells.Item(99, 13).Value = -10401.8
$ws.Cells.Item(99, 14).Value = -4497.75

# CRP!row113
$ws.Cells.Item(113, 8).Value = 1925.7858
$ws.Cells.Item(113, 9).Value = 1694.4706
$ws.Cells.Item(113, 11).Value = 1694.4706
$ws.Cells.Item(113, 13).Value = 475.5293999999999

# CRP!row126
$ws.Cells.Item(126, 8).Value = 5501
$ws.Cells.Item(126, 9).Value = 11899.8
$ws.Cells.Item(126, 10).Value = 1501.75
$ws.Cells.Item(126, 11).Value = 35699.39999999999
$ws.Cells.Item(126, 12).Value = 4505.25
$ws.Cells.Item(126, 13).Value = -33229.39999999999
$ws.Cells.Item(126, 14).Value = -9445.25

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1146.7
$ws.Cells.Item(122, 9).Value = 612.5714
$ws.Cells.Item(122, 10).Value = 2393
$ws.Cells.Item(122, 11).Value = 5513.1426
$ws.Cells.Item(122, 12).Value = 21537
$ws.Cells.Item(122, 13).Value = -3063.1426
$ws.Cells.Item(122, 14).Value = -26437

# CUL!row131
$ws.Cells.Item(131, 8).Value = 2223304.5
$ws.Cells.Item(131, 10).Value = 1275.3667
$ws.Cells.Item(131, 12).Value = 3826.1001
$ws.Cells.Item(131, 14).Value = -13906.1001

# GSM!row24
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 6225.5
$ws.Cells.Item(24, 10).Value = 6225.5
$ws.Cells.Item(24, 12).Value = 6225.5
$ws.Cells.Item(24, 14).Value = -6571.5

# GSM!row97
$ws.Cells.Item(97, 8).Value = 1683.75
$ws.Cells.Item(97, 9).Value = 1759.9
$ws.Cells.Item(97, 11).Value = 1759.9
$ws.Cells.Item(97, 13).Value = -1263.9

# GSM!row122
$ws.Cells.Item(122, 8).Value = 27305488
$ws.Cells.Item(122, 9).Value = 48402800
$ws.Cells.Item(122, 10).Value = 3082.2942
$ws.Cells.Item(122, 11).Value = 145208400
$ws.Cells.Item(122, 12).Value = 9246.882599999999
$ws.Cells.Item(122, 13).Value = -145205950
$ws.Cells.Item(122, 14).Value = -14146.8826

# GSM!row123
$ws.Cells.Item(123, 8).Value = 19699.5
$ws.Cells.Item(123, 10).Value = 20068.45
$ws.Cells.Item(123, 12).Value = 20068.45
$ws.Cells.Item(123, 14).Value = -24968.45

# GSM!row126
$ws.Cells.Item(126, 8).Value = 5439.1934
$ws.Cells.Item(126, 9).Value = 6738.6
$ws.Cells.Item(126, 10).Value = 3076.6365
$ws.Cells.Item(126, 11).Value = 20215.8
$ws.Cells.Item(126, 12).Value = 9229.9095
$ws.Cells.Item(126, 13).Value = -17745.8
$ws.Cells.Item(126, 14).Value = -14169.9095

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1925.4
$ws.Cells.Item(22, 9).Value = 322.25
$ws.Cells.Item(22, 10).Value = 2172.0386
$ws.Cells.Item(22, 11).Value = 322.25
$ws.Cells.Item(22, 12).Value = 2172.0386
$ws.Cells.Item(22, 13).Value = -27.25
$ws.Cells.Item(22, 14).Value = -2762.0386

# LTW!row27
$ws.Cells.Item(27, 8).Value = 1925.4
$ws.Cells.Item(27, 9).Value = 322.25
$ws.Cells.Item(27, 10).Value = 2172.0386
$ws.Cells.Item(27, 11).Value = 322.25
$ws.Cells.Item(27, 12).Value = 2172.0386
$ws.Cells.Item(27, 13).Value = -215.25
$ws.Cells.Item(27, 14).Value = -2386.0386

# LTW!row93
$ws.Cells.Item(93, 8).Value = 1250
$ws.Cells.Item(93, 10).Value = 1250
$ws.Cells.Item(93, 12).Value = 1250
$ws.Cells.Item(93, 14).Value = -3746

# LTW!row122
$ws.Cells.Item(122, 8).Value = 6271336
$ws.Cells.Item(122, 9).Value = 6502033.5
$ws.Cells.Item(122, 10).Value = 5002500
$ws.Cells.Item(122, 11).Value = 19506100.5
$ws.Cells.Item(122, 12).Value = 15007500
$ws.Cells.Item(122, 13).Value = -19503650.5
$ws.Cells.Item(122, 14).Value = -15012400

# LTW!row133
$ws.Cells.Item(133, 8).Value = 97965.2
$ws.Cells.Item(133, 10).Value = 97965.2
$ws.Cells.Item(133, 12).Value = 97965.2
$ws.Cells.Item(133, 14).Value = -103025.2

# WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 71429544
$ws.Cells.Item(107, 9).Value = 111111920
$ws.Cells.Item(107, 10).Value = 1280
$ws.Cells.Item(107, 11).Value = 333335760
$ws.Cells.Item(107, 12).Value = 3840
$ws.Cells.Item(107, 13).Value = -333333840
$ws.Cells.Item(107, 14).Value = -7680

# WVR!row122
$ws.Cells.Item(122, 8).Value = 1511.2222
$ws.Cells.Item(122, 9).Value = 1200.1428
$ws.Cells.Item(122, 10).Value = 2600
$ws.Cells.Item(122, 11).Value = 3600.4284
$ws.Cells.Item(122, 12).Value = 7800
$ws.Cells.Item(122, 13).Value = -1150.4284
$ws.Cells.Item(122, 14).Value = -12700

# WVR!row136
$ws.Cells.Item(136, 8).Value = 2536.7441
$ws.Cells.Item(136, 9).Value = 2743
$ws.Cells.Item(136, 10).Value = 2250.2778
$ws.Cells.Item(136, 11).Value = 8229
$ws.Cells.Item(136, 12).Value = 6750.8334
$ws.Cells.Item(136, 13).Value = -5679
$ws.Cells.Item(136, 14).Value = -11850.8334
